# ---------------------------------------------------------------------------
# Commit: "Added the Bias Ratio and updated documentation."
#
# 1. Merge "Table 2.4" and "Table 2.5" into a single sheet "Table 2.4 and 2.5"
#    - Add a header row (Portfolio / Benchmark)
#    - Column A = old "Table 2.4" data (Portfolio), Column B = old
#      "Table 2.5" data (Benchmark)
#    - Replace the old Skewness/Kurtosis (Sample + Population) mini-table
#      with a two-column (Portfolio/Benchmark) version, using the built-in
#      SKEW.P function alongside the existing _xll UDFs
#    - Delete the now-redundant "Table 2.5" sheet
# 2. Rename "Table 2.7 & 2.8" to "Tables 2.7, 2.8, and 2.9" and add a new
#    "Bias Ratio" section (Number of Std Devs / Count Above / Count Below /
#    Bias Ratio) in columns G:I
# ---------------------------------------------------------------------------

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

$fmtPct4    = "0.00%"
$fmtStat4dp = "_(* #,##0.0000_);_(* \(#,##0.0000\);_(* ""-""??_);_(@_)"
$fmtPct1    = "0.0%"
$fmtComma0  = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"
$fmtComma2  = "_(* #,##0.00_);_(* \(#,##0.00\);_(* ""-""??_);_(@_)"

# ===========================================================================
# PART 1 : "Table 2.4" + "Table 2.5"  ->  "Table 2.4 and 2.5"
# ===========================================================================

$t24 = $wb.Worksheets.Item("Table 2.4")
$t25 = $wb.Worksheets.Item("Table 2.5")

# --- grab the Benchmark (old Table 2.5) return series before we touch it ---
$benchVals = $t25.Range("A1:A36").Value2

# --- shift the Portfolio series (& blow away the old stats block) down a row
$t24.Rows.Item(1).Insert()

# Write the new header row
$t24.Range("A1").Value2 = "Portfolio"
$t24.Range("B1").Value2 = "Benchmark"
$t24.Range("E1").Value2 = "Portfolio"
$t24.Range("F1").Value2 = "Benchmark"

# Benchmark series goes into column B, rows 2:37
$t24.Range("B2:B37").Value2 = $benchVals
$t24.Range("B2:B37").NumberFormat = $fmtPct4

# Remove the old stats labels/formulas that used to live in C:D (rows 2:8
# after the insert) - the new layout below replaces them entirely.
$t24.Range("C2:D8").Clear()

# --- new stats block (columns D:F) -----------------------------------------

# Population Skewness (built-in SKEW.P)
$t24.Range("D2").Value2 = "Population Skewness"
$t24.Range("E2").Formula = "=SKEW.P(A2:A37)"
$t24.Range("F2").Formula = "=SKEW.P(B2:B37)"

# Population Skewness (UDF, array formula)
$t24.Range("D3").Value2 = "Population Skewness"
$t24.Range("E3").FormulaArray = "=_xll.Skewness_P(A2:A37)"
$t24.Range("F3").FormulaArray = "=_xll.Skewness_P(B2:B37)"

# Population Kurtosis (UDF, array formula)
$t24.Range("D5").Value2 = "Population Kurtosis"
$t24.Range("E5").FormulaArray = "=_xll.Kurtosis_P_Excess(A2:A37)"
$t24.Range("F5").FormulaArray = "=_xll.Kurtosis_P_Excess(B2:B37)"

# JB Test (manual formula off of E3/E5)
$t24.Range("D7").Value2 = "JB Test"
$t24.Range("E7").Formula = "=COUNT(A2:A37)/6*(E3^2+E5^2/4)"
$t24.Range("F7").Formula = "=COUNT(B2:B37)/6*(F3^2+F5^2/4)"

# JB Test (UDF, array formula)
$t24.Range("D8").Value2 = "JB Test"
$t24.Range("E8").FormulaArray = "=_xll.JarqueBeraTest(A2:A37)"
$t24.Range("F8").FormulaArray = "=_xll.JarqueBeraTest(B2:B37)"

# Sample Skewness (UDF, array formula)
$t24.Range("D10").Value2 = "Sample Skewness"
$t24.Range("E10").FormulaArray = "=_xll.Skewness_S(A2:A37)"
$t24.Range("F10").FormulaArray = "=_xll.Skewness_S(B2:B37)"

# Sample Kurtosis (UDF, array formula)
$t24.Range("D11").Value2 = "Sample Kurtosis"
$t24.Range("E11").FormulaArray = "=_xll.Kurtosis_S(A2:A37)"
$t24.Range("F11").FormulaArray = "=_xll.Kurtosis_S(B2:B37)"

$t24.Range("E2:F11").NumberFormat = $fmtStat4dp

# --- column widths ----------------------------------------------------------
$t24.Columns.Item("A").ColumnWidth = 8.14
$t24.Columns.Item("B").ColumnWidth = 10.71
$t24.Columns.Item("D").ColumnWidth = 19
$t24.Columns.Item("E").ColumnWidth = 12.71
$t24.Columns.Item("F").ColumnWidth = 10.71

# --- drop "Table 2.5" (its data now lives in column B of Table 2.4) --------
$t25.Delete()

# --- rename & finish up ------------------------------------------------------
$t24.Name = "Table 2.4 and 2.5"
$t24.Activate()
$t24.Range("E13").Select()

# ===========================================================================
# PART 2 : "Table 2.7 & 2.8"  ->  "Tables 2.7, 2.8, and 2.9"  (+ Bias Ratio)
# ===========================================================================

$t27 = $wb.Worksheets.Item("Table 2.7 & 2.8")

$t27.Range("H9").Value2 = "Portfolio"
$t27.Range("I9").Value2 = "Benchmark"

$t27.Range("G10").Value2 = "Number of Std Devs"
$t27.Range("H10").Value2 = 1
$t27.Range("I10").Formula = "=H10"
$t27.Range("H10:I10").NumberFormat = $fmtComma2

$t27.Range("G11").Value2 = "Count Above"
$t27.Range("H11").FormulaArray = "=SUM((A2:A37>=0)*(A2:A37<=H3))"
$t27.Range("I11").FormulaArray = "=SUM((B2:B37>=0)*(B2:B37<=I3))"

$t27.Range("G12").Value2 = "Count Below"
$t27.Range("H12").FormulaArray = "=SUM((A2:A37<0)*(A2:A37>=-H3))"
$t27.Range("I12").FormulaArray = "=SUM((B2:B37<0)*(B2:B37>=-I3))"

$t27.Range("H11:I12").NumberFormat = $fmtComma0

$t27.Range("G13").Value2 = "Bias Ratio"
$t27.Range("H13").FormulaArray = "=SUM((A2:A37>=0)*(A2:A37<=H3*H10))/(1+SUM((A2:A37<0)*(A2:A37>=-H3*H10)))"
$t27.Range("I13").FormulaArray = "=SUM((B2:B37>=0)*(B2:B37<=I3*I10))/(1+SUM((B2:B37<0)*(B2:B37>=-I3*I10)))"

$t27.Range("G14").Value2 = "Bias Ratio"
$t27.Range("H14").FormulaArray = "=_xll.BiasRatio(A2:A37,H10)"
$t27.Range("I14").FormulaArray = "=_xll.BiasRatio(B2:B37,I10)"

$t27.Range("H13:I14").NumberFormat = $fmtComma2

# --- column widths ----------------------------------------------------------
$t27.Columns.Item("G").ColumnWidth = 18.71
$t27.Columns.Item("H").ColumnWidth = 11.71
$t27.Columns.Item("I").ColumnWidth = 11.71

# --- rename & finish up -------------------------------------------------------
$t27.Name = "Tables 2.7, 2.8, and 2.9"
$t27.Activate()
$t27.Range("H14").Select()
